# Add a new "Scatter" worksheet as the last sheet (after "Line"),
# populate it with the same data used by the other chart sheets, then
# add an XY scatter chart plotting column B (x) against column A (y).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Scatter"

$values = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
    $ws.Cells.Item($row, 2).Value = [double]($values[$i] * $values[$i])
}

# xlXYScatterLines = 74 -> scatterStyle "lineMarker"
$chartObj = $ws.ChartObjects().Add(10, 200, 400, 200)
$chart = $chartObj.Chart
$chart.ChartType = 74

$chart.SeriesCollection().NewSeries()
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Scatter'!`$B`$1:`$B`$10"
$series.Values = "='Scatter'!`$A`$1:`$A`$10"

$xAxis = $chart.Axes(1)
$yAxis = $chart.Axes(2)

$xAxis.MinimumScale = 0
$xAxis.MaximumScale = 90
$yAxis.MinimumScale = 0
$yAxis.MaximumScale = 10
$yAxis.MajorUnit = 2
$xAxis.MajorUnit = 10

$chart.HasLegend = $true
$chart.Legend.Position = -4152
